# Fix the test data file (sample.05.xlsx)
#
# The real content fix: cell C20 on Sheet1 contained the string
# "list.ob.num" (a typo) which should actually read "list.obj.num"
# (consistent with the sibling "list.obj.str" value used at C16).
#
# Fixing that value means the old "list.ob.num" shared-string entry
# becomes unused and Excel will drop it from the shared strings table
# on save, shifting every subsequent shared-string index down by one
# -- which is exactly what the target diff shows for A8/B8/D8, B9,
# A18/B18, A21/B21/D21 and A22/B22. We don't need to touch those
# cells ourselves; simply correcting C20's text is enough to trigger
# that renumbering on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the typo'd field-name value in C20.
$ws.Range("C20").Value = "list.obj.num"

# The selection that was left in the sheet also moved, from D20 to
# C19:C21 with C19 as the active cell.
$ws.Range("C19:C21").Select()

# C19's formatting was reset back to the worksheet's default style
# (no explicit border/font/alignment), matching the style used by
# plain cells such as C20.
$ws.Range("C19").Style = "Normal"
